$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Segunda a Sexta")
$ws.Unprotect()

# Update the header title to include a line break (SETOR  IPM / PLANEJAMENTO)
$ws.Range("A1").Value = "SETOR  IPM`nPLANEJAMENTO"

# Update responsible-person names (merged cells G12:G18 ... L12:L18)
$ws.Range("G12").Value = "JORGE LUIZ C MELLO"
$ws.Range("H12").Value = "JOAO CARLOS M. SANTOS "
$ws.Range("I12").Value = "ANTONIO P DOS SANTOS"
$ws.Range("J12").Value = "DEVANI CORREA "
$ws.Range("K12").Value = "JOSE CARLOS VILELA"
$ws.Range("L12").Value = "LUCIANO DA SILVA SILVEIRA"

# Update registration numbers (merged cells G19:G20 ... L19:L20)
# These must be stored as TEXT (not numeric), matching the original workbook's
# shared-string representation, so force the number format to Text first.
$ws.Range("G19:L19").NumberFormat = "@"
$ws.Range("G19").Value = "1091"
$ws.Range("H19").Value = "1625"
$ws.Range("I19").Value = "2088"
$ws.Range("J19").Value = "2412"
$ws.Range("K19").Value = "2831"
$ws.Range("L19").Value = "3687"
